$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 69, shifting existing rows 69..343 down to 70..344
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new record
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "Femacal de La Calera"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = 44764
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = 100112039
$ws.Range("G69").Value = "Ciboulette"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 120
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 1500
$ws.Range("M69").Value = 1500
$ws.Range("N69").Value = "$/docena de atados"
$ws.Range("O69").Value = "Provincia de Quillota"
$ws.Range("P69").Value = 500
$ws.Range("Q69").Value = 3
$ws.Range("R69").Value = "Hortaliza"
